$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q1" sheet right before the "总计" sheet ---
$totalSheetBefore = $wb.Worksheets.Item(3)
$newQ = $wb.Worksheets.Add($totalSheetBefore)
$newQ.Name = "2022-Q1"

# Re-resolve "总计" AFTER the Add() call -- inserting a sheet shifts its
# position, and the old object reference becomes stale.
$total = $wb.Worksheets.Item(4)
$q4 = $wb.Worksheets.Item(2)

# --- 2. Populate the new "2022-Q1" sheet ---
# Column headers + A2 index cell share the same style ("s=2") already used
# on the "2021-Q4" sheet, so copy that formatting across instead of trying
# to rebuild it property-by-property.
$q4.Range("B1:H1").Copy($newQ.Range("B1:H1"))
$q4.Range("A2").Copy($newQ.Range("A2"))

$newQ.Range("A2").Value = 0

$newQ.Range("B2:G2").NumberFormat = "@"
$newQ.Range("B2").Value = "320017"
$newQ.Range("C2").Value = "诺安全球收益不动产(QDII)"
$newQ.Range("D2").Value = "0.29"
$newQ.Range("E2").Value = "93.32"
$newQ.Range("F2").Value = "4.86"
$newQ.Range("G2").Value = "0.0141"
$newQ.Range("B2:G2").ClearFormats()

$newQ.Range("H2").Value = 10

# --- 3. Update the "总计" roll-up sheet with the new quarter on top ---
$total.Rows.Item(2).Insert()

# Carry the index-column style ("s=2") down into the freshly inserted row.
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
